$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text (matches the
# workbook's original inline-string cells) so that number-looking
# values such as "482.92" are not re-interpreted as numeric cells.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "68.107.03"
$ws.Range("D3").Value = "3.884.30"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "482.92"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").Value = "144.64"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.723"
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +4.06%  "
$ws.Range("D11").Value = "0.0000353"
$ws.Range("E11").Value = "  +12.59%  "
$ws.Range("D12").Value = "42.65"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").Value = "10.65"
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").Value = "4.504.24"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "14.57"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").Value = "3.882.37"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "19.69"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").Value = "68.174.71"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").Value = "434.14"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "3.38"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "14.64"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").Value = "88.02"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").Value = "11.49"
$ws.Range("E25").Value = "  +15.75%  "
$ws.Range("D26").Value = "3.57"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").Value = "10.49"
$ws.Range("E27").Value = "  +5.43%  "
$ws.Range("D28").Value = "37.97"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "5.80"
$ws.Range("E29").Value = "  +4.69%  "
$ws.Range("D30").Value = "700.62"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").Value = "13.36"
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("E34").Value = "  +33.62%  "
$ws.Range("D35").Value = "41.48"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").Value = "59.62"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("D37").Value = "5.71"
$ws.Range("E37").Value = "  +3.73%  "
$ws.Range("E38").Value = "  -7.05%  "
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "0.0473"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "3.03"
$ws.Range("E41").Value = "  +3.26%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "2.73"
$ws.Range("E42").Value = "  +7.17%  "
$ws.Range("D43").Value = "3.01"
$ws.Range("E43").Value = "  +8.82%  "
$ws.Range("D44").Value = "0.340"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").Value = "146.03"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("D50").Value = "3.12"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").Value = "2.83"
$ws.Range("E51").Value = "  -2.65%  "

# Restore the original (default) cell style now that the text values
# are locked in, so no stray style indices are left on the cells.
$dataRange.Style = "Normal"
